$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The grid's "ok" marker block in column D shifts down by 11 rows:
# it used to flag D66:D77, now it flags D77:D88.
$xlPasteFormats = -4122

# 1) D78:D88 become new "ok" cells. They already carry the right look, so a
#    plain value write is enough.
$ws.Range("D78:D88").Value = "ok"

# 2) D77 stays "ok" but becomes the top/highlighted row of the block, so give
#    it the formatting the old top row (D66) used to have.
$ws.Range("D66").Copy()
$ws.Range("D77").PasteSpecial($xlPasteFormats)
$ws.Range("D77").Value = "ok"

# 3) D66:D76 are no longer part of the "ok" block: clear their content and
#    restore the same "cleared" formatting already used elsewhere in the
#    sheet for non-flagged cells (e.g. D46:D57).
$ws.Range("D46").Copy()
$ws.Range("D66:D76").PasteSpecial($xlPasteFormats)
$ws.Range("D66:D76").ClearContents()
